## httk-benchmarks.xlsx update: add the "2.7.0" release row to the
## benchmark table on Sheet1 (table grows from A1:R29 to A1:R30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Grow the Excel Table (ListObject) by one row -------------------------
$tbl = $ws.ListObjects.Item("Table1")
$lastRowNum = $tbl.Range.Row + $tbl.Range.Rows.Count - 1   # 29
$newRowNum  = $lastRowNum + 1                              # 30

$newListRow = $tbl.ListRows.Add()

# Copy the formatting (style) of the previous last row onto the new row so
# the new cells pick up the same left-aligned style (s="1") the rest of the
# table uses.
$srcRange = $ws.Range("A" + $lastRowNum + ":R" + $lastRowNum)
$dstRange = $ws.Range("A" + $newRowNum + ":R" + $newRowNum)
$srcRange.Copy($dstRange)

# --- Fill in the new row's values ------------------------------------------
$row = $newRowNum

$ws.Range("A$row").Value = "2.7.0"
$ws.Range("B$row").Value = 1026
$ws.Range("C$row").Value = 0.99980000000000002
$ws.Range("D$row").Value = 1
$ws.Range("E$row").Value = 0.99990000000000001
$ws.Range("F$row").Value = 0.93500000000000005
$ws.Range("G$row").Value = 352
$ws.Range("H$row").Value = 0.2712
$ws.Range("I$row").Value = 352
$ws.Range("J$row").Value = 2.3740000000000001
$ws.Range("K$row").Value = 43
$ws.Range("L$row").Value = 1.5309999999999999
$ws.Range("M$row").Value = 160
$ws.Range("N$row").Value = 1.202
$ws.Range("O$row").Value = 160
$ws.Range("P$row").Value = 0.629
$ws.Range("Q$row").Value = 863
$ws.Range("R$row").Value = "Updated CvT data, added pfas, dermal models, updated IVD models"

# --- Misc view bookkeeping (best effort, mirrors the authored file) -------
$ws.Activate()
$excel.ActiveWindow.Zoom = 145
$ws.Range("J$row").Select()
